$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CottonObserved")

# Remove the two observed-data columns "Cotton.Seed.NConc" (BR) and
# "Cotton.Seed.N" (BS). Deleting the whole columns shifts every cell,
# formula, column-width definition, row/dimension reference etc. two
# columns to the left automatically - exactly like doing it by hand in
# Excel (select columns BR:BS > right-click > Delete).
$ws.Range("BR1:BS1").EntireColumn.Delete()

# The AutoFilter's stored range doesn't auto-shrink with the column
# delete, so refresh it over the worksheet's new used range.
$ws.AutoFilterMode = $false
$ws.Range("A1:ET44").AutoFilter()

# Likewise, the workbook-level hidden "_FilterDatabase" name still
# points at the old (wider) range; update it to match.
foreach ($n in $wb.Names) {
    if ($n.Name -eq "CottonObserved!_FilterDatabase") {
        $n.RefersTo = "=CottonObserved!`$A`$1:`$ET`$44"
    }
}

# Restore the user's on-screen selection to where the deleted columns
# used to be.
$ws.Activate()
$ws.Range("BR4").Select()
